$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.532.94"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "3.113.02"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "525.95"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.11"
$ws.Range("E6").Value = "  -2.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.111.09"
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.25"
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("E12").Value = "  +3.23%  "
$ws.Range("D13").Value = "3.647.53"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("E14").Value = "  +3.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.35"
$ws.Range("E15").Value = "  -2.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000164"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "57.632.48"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").Value = "3.109.73"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("E19").Value = "  -2.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.46"
$ws.Range("E20").Value = "  -2.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.91"
$ws.Range("E21").Value = "  -1.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "347.93"
$ws.Range("E22").Value = "  +2.88%  "
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.07"
$ws.Range("E25").Value = "  +2.21%  "
$ws.Range("E26").Value = "  -1.85%  "
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.43"
$ws.Range("E30").Value = "  +3.48%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +0.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.02"
$ws.Range("E33").Value = "  -7.28%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.98"
$ws.Range("E35").Value = "  +7.57%  "
$ws.Range("E36").Value = "  -2.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.04"
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.96"
$ws.Range("E39").Value = "  -4.81%  "
$ws.Range("E40").Value = "  -3.18%  "
$ws.Range("E41").Value = "  +6.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0666"
$ws.Range("E42").Value = "  +1.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.62"
$ws.Range("E43").Value = "  +6.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.701"
$ws.Range("E44").Value = "  +2.31%  "
$ws.Range("D45").Value = "3.150.31"
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("D46").Value = "2.355.69"
$ws.Range("E46").Value = "  +1.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.51"
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("E49").Value = "  +3.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.961"
$ws.Range("E50").Value = "  -1.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.03"
$ws.Range("E51").Value = "  +0.19%  "
